# Insert a new data row at row 18 (pushing existing rows 18..87 down to 19..88)
# and populate it with the new "Femacal de La Calera - Arveja Verde" record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()

$ws.Cells.Item(18, 1).Value  = 3
$ws.Cells.Item(18, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(18, 3).Value  = "Coquimbo"
$ws.Cells.Item(18, 4).Value  = 45114
$ws.Cells.Item(18, 5).Value  = 5
$ws.Cells.Item(18, 6).Value  = 100112022
$ws.Cells.Item(18, 7).Value  = "Arveja Verde"
$ws.Cells.Item(18, 8).Value  = "Perfection"
$ws.Cells.Item(18, 9).Value  = "Primera"
$ws.Cells.Item(18, 10).Value = 35
$ws.Cells.Item(18, 11).Value = 23000
$ws.Cells.Item(18, 12).Value = 23000
$ws.Cells.Item(18, 13).Value = 23000
$ws.Cells.Item(18, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 920
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
